# Cleaned up the data processing: insert new "most_frequent_count" and "csim"
# metadata columns (computed when generating dataset metadata), shifting the
# existing memory_consumed_bytes / pattern_count / patterns columns right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the old "memory_consumed_bytes" column (T),
# pushing T,U,V -> V,W,X and leaving T,U empty for the new fields.
$ws.Columns("T:U").Insert()

# New header cells
$ws.Range("T1").Value = "most_frequent_count"
$ws.Range("U1").Value = "csim"

# Per-attribute "position" (B) and "type" (C) columns were reworked.
$ws.Range("B2").Value = -1
$ws.Range("C2").Value = "date"

$ws.Range("B3").Value = -1
$ws.Range("C3").Value = "date"

$ws.Range("B4").Value = -1
$ws.Range("C4").Value = "date"

$ws.Range("B5").Value = -1
$ws.Range("C5").Value = "string"

$ws.Range("B6").Value = -1
$ws.Range("C6").Value = "string"

$ws.Range("B7").Value = -1
$ws.Range("C7").Value = "date"

$ws.Range("B8").Value = -1
$ws.Range("C8").Value = "date"

$ws.Range("B9").Value = -1
$ws.Range("C9").Value = "date"

$ws.Range("B10").Value = -1
$ws.Range("C10").Value = "date"

# Row 2 ("Licence number") statistics were recomputed.
$ws.Range("E2").Value = 4604
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = -1
$ws.Range("L2").Value = -1

# New per-row most_frequent_count / csim values.
$ws.Range("T2").Value = 4
$ws.Range("U2").Value = 0

$ws.Range("T3").Value = 4033
$ws.Range("U3").Value = 0

$ws.Range("T4").Value = 4494
$ws.Range("U4").Value = 0

$ws.Range("T5").Value = 144
$ws.Range("U5").Value = 1

$ws.Range("T6").Value = 145
$ws.Range("U6").Value = 1

$ws.Range("T7").Value = 1185
$ws.Range("U7").Value = 0

$ws.Range("T8").Value = 4
$ws.Range("U8").Value = 0

$ws.Range("T9").Value = 976
$ws.Range("U9").Value = 0

$ws.Range("T10").Value = 786
$ws.Range("U10").Value = 0
